# Applies the authoring edit described by the commit:
#  1. Removes the second slide ("Luís Novaes" persona slide, sldId 256).
#  2. Refreshes the cached text of the auto-updating Date and Slide-Number
#     placeholder fields on the slide master and every slide layout
#     (12/03/2020 -> 01/05/2020, and the slide-number preview glyph
#     ‹#› -> ‹nº›).

$p = $ppt.ActivePresentation

# --- 1. Delete the "Luís Novaes" slide (second slide in the deck) ---------
$p.Slides.Item(2).Delete()

# --- 2. Refresh date / slide-number placeholder text on master + layouts --
$newDate = "01/05/2020"
$newSlideNum = [string][char]0x2039 + "n" + [string][char]0xBA + [string][char]0x203A

function Update-HeaderFooterPlaceholders($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        $placeholderType = -1
        try {
            $placeholderType = $shp.PlaceholderFormat.Type
        } catch {
            $placeholderType = -1
        }

        if ($placeholderType -eq 16) {
            # ppPlaceholderDate
            $shp.TextFrame.TextRange.Text = $newDate
        } elseif ($placeholderType -eq 13) {
            # ppPlaceholderSlideNumber
            $shp.TextFrame.TextRange.Text = $newSlideNum
        }
    }
}

$master = $p.SlideMaster
Update-HeaderFooterPlaceholders $master

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-HeaderFooterPlaceholders $master.CustomLayouts.Item($li)
}
